# "Fix back image bug"
# - Remove the second slide (the duplicate/"back" slide, id 257) from the deck.
# - Replace the lyrics placeholder text on the remaining slide with "Prequel text".

$p = $ppt.ActivePresentation

# Slide 1 still holds all the lyrics in its single "Rectangle 1" shape;
# collapse that whole text body down to a single run reading "Prequel text".
$s1 = $p.Slides.Item(1)
$shp = $s1.Shapes.Item(1)
$shp.TextFrame.TextRange.Text = "Prequel text"

# Slide 2 (sldId 257) is dropped entirely from the presentation.
$p.Slides.Item(2).Delete()
